$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final dataset for the "Deudores" table (rows 2..30), sorted by Cliente,
# reflecting new debtors added, one removed (CARNIVOROS), and several
# Fecha/Valor updates.
$data = @(
  @(1, "ALISO", 46039, 50000, $false),
  @(2, "ARROZ PAISA SUBA", 46043, 166000, $false),
  @(3, "CALDAS WOK", 46035, 85000, $false),
  @(4, "CAMPO VERDE TOCANCIPA", 46044, 285000, $false),
  @(5, "CAMPO VERDE ZIPAQUIRA", 46044, 483000, $false),
  @(6, "CARNES JOHANA", 46041, 164000, $false),
  @(7, "CARNILANDIA", 46042, 170000, $false),
  @(8, "CIMARRON DORADO", 46041, 426500, $false),
  @(9, "COCINA CHINA", 46043, 170000, $false),
  @(10, "COCINA CHINA", 46039, 170000, $false),
  @(11, "CRISTIAN ACACIAS", 46042, 1000000, $false),
  @(12, "DARWIN FUTBOL", 45921, 200000, $false),
  @(13, "DAVIDCITO", 45947, 100000, $false),
  @(14, "DOÑA SANDRA", 46039, 100000, $false),
  @(15, "FRANCO", 45996, 20000, $false),
  @(16, "JULIANA POLLO", 46042, 264000, $false),
  @(17, "LA 13", 46044, 299000, $false),
  @(18, "LA SELECTA", 45912, 82000, $false),
  @(19, "LOS PAISANOS", 46044, 453000, $false),
  @(20, "MERKA FRUVER DEXI", 45995, 339000, $false),
  @(21, "MERKA FRUVER DEXI", 45988, 15400, $false),
  @(22, "MULTICARNES", 46038, 646900, $false),
  @(23, "NEVADA", 46038, 175800, $false),
  @(24, "NEVADA", 46031, 21900, $false),
  @(25, "PARAISO FUNZA", 46038, 76000, $false),
  @(26, "PINILLA", 45931, 82000, $false),
  @(27, "PINILLA SOACHA", 46039, 266000, $false),
  @(28, "PLAZA JESSICA", 46039, 1349200, $false),
  @(29, "SANTANDER SUR", 46042, 176000, $false)
)

for ($i = 0; $i -lt $data.Length; $i++) {
  $r = $i + 2
  $row = $data[$i]
  $ws.Cells.Item($r, 1).Value = $row[0]
  $ws.Cells.Item($r, 2).Value = $row[1]
  $ws.Cells.Item($r, 3).Value = $row[2]
  $ws.Cells.Item($r, 3).NumberFormat = "YYYY-MM-DD"
  $ws.Cells.Item($r, 4).Value = $row[3]
  $ws.Cells.Item($r, 5).Value = $row[4]
}
